$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCRSource")

# Remove the "forward_primer" and "reverse_primer" columns (A and B),
# shifting the remaining columns left by two.
$ws.Columns.Item(1).Delete()
$ws.Columns.Item(1).Delete()
